# expandVarAll / expandVar support data:
#  - "actions" sheet gets a new "expand" action row for the " Vaccs" variable,
#    mapping pfizer/moderna codes via an R list() passed as the action's parameters.
#  - "dicos" sheet gets a new YESNO dictionary (YES=1 / NO=2) used by that action.
#  - "dicos" becomes the active/selected sheet (was "dictionary").

$wb = $excel.ActiveWorkbook
$wsDict = $wb.Worksheets.Item("dictionary")
$wsDicos = $wb.Worksheets.Item("dicos")
$wsActions = $wb.Worksheets.Item("actions")

# --- actions sheet: header gains "action_group", data row for " Vaccs" -> expand ---
# (write order chosen so new shared-string indices come out expand, " Vaccs",
#  action_group, list(...) - matching how the workbook was actually authored)
$wsActions.Range("B2").Value = "expand"
$wsActions.Range("A2").Value = " Vaccs"
$wsActions.Range("B1").Value = "action_group"
$wsActions.Range("C2").Value = 'list(pfizer="pfizer", "moderna"="moderna")'

$wsActions.Range("B1").Font.Bold = $true
$wsActions.Range("A2:C2").Font.Bold = $true

# extra (empty) formatted row below the data
$wsActions.Range("C3").Font.Bold = $true

# column widths for the now-wider "actions" content
$wsActions.Columns.Item(1).ColumnWidth = 15.498697916666666
$wsActions.Columns.Item(2).ColumnWidth = 12.998697916666666
$wsActions.Columns.Item(3).ColumnWidth = 54.498697916666664

$wsActions.Activate()
$wsActions.Range("C3").Select()
$excel.ActiveWindow.Zoom = 200

# --- dicos sheet: YESNO dictionary (YES=1, NO=2) ---
$wsDicos.Range("A2").Value = "YESNO"
$wsDicos.Range("B2").Value = "YES"
$wsDicos.Range("C2").Value = 1
$wsDicos.Range("A3").Value = "YESNO"
$wsDicos.Range("B3").Value = "NO"
$wsDicos.Range("C3").Value = 2

$wsDicos.Range("A2:B3").Font.Bold = $true

# dicos becomes the active / selected sheet
$wsDicos.Activate()
$wsDicos.Range("C4").Select()

Write-Output "edit applied"
